$d = $word.ActiveDocument

# First paragraph: "This is a Microsoft word document."
$p1 = $d.Paragraphs(1)

# Add two trailing spaces to the existing sentence (excluding the
# paragraph mark from the range we touch).
$body = $p1.Range
$body.MoveEnd(1, -1) | Out-Null
$body.InsertAfter("  ")

# Append "(This is a change – Version for main branch)" as three
# separate red-colored runs, mirroring the source edit.
$run2 = $p1.Range
$run2.MoveEnd(1, -1) | Out-Null
$run2.Collapse(0) | Out-Null
$run2.InsertAfter("(This is a change – Ve")
$run2.Font.Color = 255

$run3 = $p1.Range
$run3.MoveEnd(1, -1) | Out-Null
$run3.Collapse(0) | Out-Null
$run3.InsertAfter("rsion for main branch")
$run3.Font.Color = 255

$run4 = $p1.Range
$run4.MoveEnd(1, -1) | Out-Null
$run4.Collapse(0) | Out-Null
$run4.InsertAfter(")")
$run4.Font.Color = 255
